# Weekly fruit/vegetable price update: insert a new daily record as row 596
# (shifting all subsequent rows down by one) on the "Zapallo italiano" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 596, pushing existing rows 596-645 down
# to 597-646 (dimension grows from R645 to R646).
$ws.Rows.Item(596).Insert()

# Populate the newly inserted row with the new market observation.
$ws.Range("A596").Value = 6
$ws.Range("B596").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C596").Value = "Metropolitana"
$ws.Range("D596").Value = 45223
$ws.Range("E596").Value = 13
$ws.Range("F596").Value = 100112032
$ws.Range("G596").Value = "Zapallo italiano"
$ws.Range("H596").Value = "Sin especificar"
$ws.Range("I596").Value = "Primera"
$ws.Range("J596").Value = 320
$ws.Range("K596").Value = 13000
$ws.Range("L596").Value = 13000
$ws.Range("M596").Value = 13000
$ws.Range("N596").Value = "`$/caja 50 unidades"
$ws.Range("O596").Value = "Región de O'Higgins"
$ws.Range("P596").Value = 260
$ws.Range("Q596").Value = 50
$ws.Range("R596").Value = "Hortaliza"
